# "add rural communities interactions back in, fix de_dg files (not run yet)"
#
# The stakeholder-interaction matrix had a stray "Legislature" row (row 11)
# that no longer belongs in the list. Delete the entire row: Excel shifts
# every row below it up by one (so "Friant Water Authority" etc. move into
# place) and, because "Legislature" is no longer referenced anywhere, the
# shared-string table drops it on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows("11:11").Delete()

# Leave the selection on the row that used to be row 12 (now row 11), same
# as where the edit was made.
$ws.Range("A11:XFD11").Select()
